$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.571.66"
$ws.Range("E2").Value = "'  -0.36%  "
$ws.Range("D3").Value = "'1.874.58"
$ws.Range("E3").Value = "'  -1.12%  "
$ws.Range("E4").Value = "'  +0.12%  "
$ws.Range("D5").Value = "'247.73"
$ws.Range("E5").Value = "'  +0.90%  "
$ws.Range("D6").Value = "'1.0000"
$ws.Range("E6").Value = "'  +0.07%  "
$ws.Range("D7").Value = "'0.4736"
$ws.Range("E7").Value = "'  -1.08%  "
$ws.Range("D8").Value = "'0.2899"
$ws.Range("E8").Value = "'  -0.21%  "
$ws.Range("D9").Value = "'0.06475"
$ws.Range("E9").Value = "'  -1.50%  "
$ws.Range("D10").Value = "'22.00"
$ws.Range("E10").Value = "'  +2.05%  "
$ws.Range("D11").Value = "'0.07729"
$ws.Range("E11").Value = "'  -0.75%  "
$ws.Range("D12").Value = "'0.7417"
$ws.Range("E12").Value = "'  -0.09%  "
$ws.Range("D13").Value = "'1.875.07"
$ws.Range("E13").Value = "'  -1.13%  "
$ws.Range("D14").Value = "'96.02"
$ws.Range("E14").Value = "'  -0.87%  "
$ws.Range("D15").Value = "'5.166"
$ws.Range("E15").Value = "'  -0.53%  "
$ws.Range("D16").Value = "'274.44"
$ws.Range("E16").Value = "'  -3.10%  "
$ws.Range("D17").Value = "'30.643.13"
$ws.Range("E17").Value = "'  -0.01%  "
$ws.Range("D18").Value = "'13.27"
$ws.Range("E18").Value = "'  -3.16%  "
$ws.Range("E19").Value = "'  +0.00%  "
$ws.Range("D20").Value = "'0.000007468"
$ws.Range("E20").Value = "'  -2.27%  "
$ws.Range("D21").Value = "'2.120.98"
$ws.Range("E21").Value = "'  -0.84%  "
$ws.Range("D22").Value = "'1.000"
$ws.Range("E22").Value = "'  +0.19%  "
$ws.Range("D23").Value = "'5.199"
$ws.Range("E23").Value = "'  -2.39%  "
$ws.Range("D24").Value = "'6.171"
$ws.Range("E24").Value = "'  -1.46%  "
$ws.Range("D25").Value = "'9.199"
$ws.Range("E25").Value = "'  -1.62%  "
$ws.Range("D26").Value = "'164.92"
$ws.Range("E26").Value = "'  -0.51%  "
$ws.Range("D27").Value = "'18.67"
$ws.Range("E27").Value = "'  -2.57%  "
$ws.Range("E28").Value = "'  -5.15%  "
$ws.Range("D29").Value = "'0.09938"
$ws.Range("E29").Value = "'  -0.60%  "
$ws.Range("D30").Value = "'1.346"
$ws.Range("E30").Value = "'  -2.66%  "
$ws.Range("D31").Value = "'1.509"
$ws.Range("E31").Value = "'  -0.81%  "
$ws.Range("D32").Value = "'4.238"
$ws.Range("E32").Value = "'  -3.08%  "
$ws.Range("E33").Value = "'  -1.58%  "
$ws.Range("D34").Value = "'0.04769"
$ws.Range("E34").Value = "'  -0.42%  "
$ws.Range("D35").Value = "'1.119"
$ws.Range("E35").Value = "'  -1.68%  "
$ws.Range("D36").Value = "'0.6917"
$ws.Range("E36").Value = "'  -2.32%  "
$ws.Range("D37").Value = "'2.718"
$ws.Range("E37").Value = "'  -0.01%  "
$ws.Range("D38").Value = "'0.01847"
$ws.Range("E38").Value = "'  -1.75%  "
$ws.Range("D39").Value = "'2.752"
$ws.Range("E39").Value = "'  -0.75%  "
$ws.Range("E40").Value = "'  -4.45%  "
$ws.Range("D41").Value = "'73.12"
$ws.Range("E41").Value = "'  +2.59%  "
$ws.Range("D42").Value = "'1.965"
$ws.Range("E42").Value = "'  +1.43%  "
$ws.Range("E43").Value = "'  +0.10%  "
$ws.Range("E44").Value = "'  -1.28%  "
$ws.Range("D45").Value = "'0.8334"
$ws.Range("E45").Value = "'  -2.29%  "
$ws.Range("D46").Value = "'101.30"
$ws.Range("E46").Value = "'  -1.60%  "
$ws.Range("D47").Value = "'9.353"
$ws.Range("E47").Value = "'  -1.52%  "
$ws.Range("D48").Value = "'35.34"
$ws.Range("E48").Value = "'  -0.37%  "
$ws.Range("E49").Value = "'  -2.99%  "
$ws.Range("D50").Value = "'913.86"
$ws.Range("E50").Value = "'  -2.42%  "
$ws.Range("E51").Value = "'  +0.88%  "
